# Add a new "level up" reward row (row 11) to the LevelInfo table,
# mirroring the style/format of the preceding data row (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 10) down into the
# new row 11 so the new row matches the table's existing style.
$ws.Range("A10:E10").Copy()
[void]$ws.Range("A11:E11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values.
$ws.Range("A11").Value = 100
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = "|您获得了一个|Yellow|经典卡包||。"
$ws.Range("E11").Value = 1

# Grow the table (表1) to include the new row.
$tbl = $ws.ListObjects.Item(1)
[void]$tbl.Resize($ws.Range("A3:E11"))

# Match the resulting selection state.
[void]$ws.Range("E11").Select()
